# Replace the placeholder paragraph with the full multi-paragraph post body.
$d = $word.ActiveDocument

$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">The reason that there has been a large influx of social scientists going into the tech industry has a few factors. </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">First, technology has is no longer just hardware and code buried deep in some server farm somewhere in Nebraska anymore. Technology has infiltrated our society at </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>it</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> basest level</w:t></w:r><w:r><w:t xml:space="preserve"> and social scientists brins a unique set of skills to the table.</w:t></w:r><w:r><w:t xml:space="preserve"> We need social scientists in the tech sector to help developers understand big questions such as the implications of the use of artificial intelligence, and how the use of social media will impact the mental health of people from childhood to adulthood. They can also help tech companies understand and navigate human behavior and social structures to facilitate the navigating the balance between technology and Society. </w:t></w:r></w:p><w:p><w:r><w:t>Second is that the tech industry is starting to recognize the importance of diversity and inclusion.</w:t></w:r><w:r><w:t xml:space="preserve"> Because of the social </w:t></w:r><w:r><w:t>scientist’s</w:t></w:r><w:r><w:t xml:space="preserve"> expertise in working with social dynamics and human interactions their insight could contribute greatly to diverse and inclusive environments by helping to address things like biases in algorithms and helping to design more inclusive products so that tech company work environments are better for both employees and the users. </w:t></w:r></w:p><w:p><w:r><w:t>The role of sociologists and social scientists in the tech industry should be to actively participate in the development process. By doing this</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space="preserve"> social scientists can be a part of creating technologies that better align with the values that their society </w:t></w:r><w:r><w:t>holds</w:t></w:r><w:r><w:t xml:space="preserve"> as important. Sociologists and social scientists also could serve as the </w:t></w:r><w:r><w:t>middleman</w:t></w:r><w:r><w:t xml:space="preserve"> between the tech companies and the </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>general public</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> by helping to translate complex technical problems into easier to understand terms and concepts suitable for consumption by the general society. </w:t></w:r></w:p><w:p><w:r><w:t>In short, the increased presence of social scientists in the tech industry just emphasizes the awareness that technological development and innovation is something that requires multiple disciplines to be executed correctly, but also in an ethical manner.</w:t></w:r></w:p>
'@

$d.Content.InsertXML($xml)
